# replace grouptc-hash-v2 to grouptc-cuckoo
# Update the "GroupTC-HS" timing column (E) and its derived speedup column (I)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.793
$ws.Range("E3").Value = 11.249
$ws.Range("E4").Value = 32.499
$ws.Range("E5").Value = 90.252
$ws.Range("E6").Value = 246.985
$ws.Range("E7").Value = 653.174
$ws.Range("E8").Value = 1710.49
$ws.Range("E9").Value = 4427.437

$ws.Range("I2").Value = 3.035855523332454
$ws.Range("I3").Value = 3.177882478442529
$ws.Range("I4").Value = 3.264100433859503
$ws.Range("I5").Value = 3.339981385454062
$ws.Range("I6").Value = 3.335793671680467
$ws.Range("I7").Value = 3.370336235061408
$ws.Range("I8").Value = 3.376357651900917
$ws.Range("I9").Value = 3.38862280818451
